$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 36
$ws.Cells.Item($row, 1).Value = "2025-04-28 23:56:01"
$ws.Cells.Item($row, 2).Value = 156
